$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "78L05G"
$ws.Range("C5").Value = "UC_SOT89"

$ws.Range("A13").Select()
